$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 416, shifting existing rows (416 onward) down by one.
$ws.Rows.Item(416).Insert()

# Populate the newly inserted row 416 with the new translation key/value.
$ws.Range("A416").Value = "errors.CROSSPLOT_SETTINGS_MISMATCH"
$ws.Range("B416").Value = "Crossplot settings mismatch. Please generate the crossplot with your most recent changes."
